# Remove the "Dividends" (H) and "Stock Splits" (I) columns from every
# ticker sheet, keeping only Date/Ticker/Open/High/Low/Close/Volume.
# Also refresh the last row (row 22) of each sheet with updated
# Close/Volume (and a couple of other slightly-revised) figures.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("H1:I1").EntireColumn.Delete()
}

$goog = $wb.Worksheets.Item("GOOG")
$goog.Cells.Item(22, 6).Value = 168.8101043701172   # F22 Close
$goog.Cells.Item(22, 7).Value = 13065359            # G22 Volume

$aapl = $wb.Worksheets.Item("AAPL")
$aapl.Cells.Item(22, 6).Value = 226.0998992919922   # F22 Close
$aapl.Cells.Item(22, 7).Value = 43997155            # G22 Volume

$amzn = $wb.Worksheets.Item("AMZN")
$amzn.Cells.Item(22, 4).Value = 185.9900054931641   # D22 Low
$amzn.Cells.Item(22, 6).Value = 185.8699951171875   # F22 Close
$amzn.Cells.Item(22, 7).Value = 25308146            # G22 Volume

$msft = $wb.Worksheets.Item("MSFT")
$msft.Cells.Item(22, 5).Value = 418.8250122070312   # E22 Low
$msft.Cells.Item(22, 6).Value = 422.5599975585938   # F22 Close
$msft.Cells.Item(22, 7).Value = 11694962            # G22 Volume
